$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4:10").Delete() | Out-Null

# Update row 2 (sale record #288)
$ws.Range("A2").Value = 288
$ws.Range("B2").Value = "Cocacola"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 8330
$ws.Range("E2").Value = "efectivo"
$ws.Range("F2").Value = "2025-08-07 14:59:26"

# Update row 3 (sale record #289)
$ws.Range("A3").Value = 289
$ws.Range("B3").Value = "Cocacola"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 8330
$ws.Range("E3").Value = "efectivo"
$ws.Range("F3").Value = "2025-08-07 15:00:50"

# Update the summary section
$ws.Range("B6").Value = 16660
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 2
